$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2469.9
$ws.Range("I19").Value = 900
$ws.Range("K19").Value = 900
$ws.Range("M19").Value = -725
$ws.Range("H46").Value = 5574.9375
$ws.Range("I46").Value = 4877
$ws.Range("J46").Value = 6272.875
$ws.Range("K46").Value = 14631
$ws.Range("L46").Value = 18818.625
$ws.Range("M46").Value = -14512
$ws.Range("N46").Value = -19056.625
$ws.Range("H60").Value = 5574.9375
$ws.Range("I60").Value = 4877
$ws.Range("J60").Value = 6272.875
$ws.Range("K60").Value = 14631
$ws.Range("L60").Value = 18818.625
$ws.Range("M60").Value = -14147
$ws.Range("N60").Value = -19786.625
$ws.Range("H80").Value = 1741.0667
$ws.Range("I80").Value = 351.7
$ws.Range("K80").Value = 1055.1
$ws.Range("M80").Value = -57.09999999999991
$ws.Range("H83").Value = 1741.0667
$ws.Range("I83").Value = 351.7
$ws.Range("K83").Value = 3165.3
$ws.Range("M83").Value = 1826.7
$ws.Range("H138").Value = 2533.1128
$ws.Range("J138").Value = 3752.1936
$ws.Range("L138").Value = 11256.5808
$ws.Range("N138").Value = -21536.5808

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2796680.8
$ws.Range("I32").Value = 3194974
$ws.Range("K32").Value = 3194974
$ws.Range("M32").Value = -3194687
$ws.Range("H74").Value = 359941.3
$ws.Range("I74").Value = 501480.6
$ws.Range("K74").Value = 501480.6
$ws.Range("M74").Value = -500606.6
$ws.Range("H77").Value = 359941.3
$ws.Range("I77").Value = 501480.6
$ws.Range("K77").Value = 2507403
$ws.Range("M77").Value = -2503035
$ws.Range("H97").Value = 884821.4399999999
$ws.Range("I97").Value = 1280391.1
$ws.Range("K97").Value = 1280391.1
$ws.Range("M97").Value = -1279895.1
$ws.Range("H122").Value = 2057.3635
$ws.Range("I122").Value = 1756.5294
$ws.Range("K122").Value = 5269.5882
$ws.Range("M122").Value = -2819.5882

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3590.4583
$ws.Range("I20").Value = 3331.5
$ws.Range("J20").Value = 4108.375
$ws.Range("K20").Value = 3331.5
$ws.Range("L20").Value = 4108.375
$ws.Range("M20").Value = -3084.5
$ws.Range("N20").Value = -4602.375
$ws.Range("H105").Value = 66687596
$ws.Range("I105").Value = 90935470
$ws.Range("K105").Value = 90935470
$ws.Range("M105").Value = -90933723
$ws.Range("H134").Value = 3859.6956
$ws.Range("I134").Value = 3262.4092
$ws.Range("K134").Value = 9787.2276
$ws.Range("M134").Value = -7252.2276

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45459424
$ws.Range("I31").Value = 90910780
$ws.Range("J31").Value = 8063.4546
$ws.Range("K31").Value = 90910780
$ws.Range("L31").Value = 8063.4546
$ws.Range("M31").Value = -90910485
$ws.Range("N31").Value = -8653.454600000001
$ws.Range("H34").Value = 45459424
$ws.Range("I34").Value = 90910780
$ws.Range("J34").Value = 8063.4546
$ws.Range("K34").Value = 90910780
$ws.Range("L34").Value = 8063.4546
$ws.Range("M34").Value = -90910578
$ws.Range("N34").Value = -8467.454600000001
$ws.Range("H62").Value = 22623.95
$ws.Range("I62").Value = 16776.223
$ws.Range("K62").Value = 16776.223
$ws.Range("M62").Value = -16152.223
$ws.Range("H65").Value = 22623.95
$ws.Range("I65").Value = 16776.223
$ws.Range("K65").Value = 83881.11500000001
$ws.Range("M65").Value = -80761.11500000001
$ws.Range("H93").Value = 10901.333
$ws.Range("J93").Value = 39995.5
$ws.Range("L93").Value = 39995.5
$ws.Range("N93").Value = -43739.5
$ws.Range("H94").Value = 6073.3335
$ws.Range("I94").Value = 9858.200000000001
$ws.Range("J94").Value = 1342.25
$ws.Range("K94").Value = 9858.200000000001
$ws.Range("L94").Value = 1342.25
$ws.Range("M94").Value = -9407.200000000001
$ws.Range("N94").Value = -2244.25
$ws.Range("H132").Value = 2820.0952
$ws.Range("I132").Value = 2043.2424
$ws.Range("K132").Value = 6129.7272
$ws.Range("M132").Value = -3599.7272

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4701.4287
$ws.Range("J39").Value = 6311
$ws.Range("L39").Value = 18933
$ws.Range("N39").Value = -19521
$ws.Range("H62").Value = 1981.9474
$ws.Range("I62").Value = 999
$ws.Range("K62").Value = 2997
$ws.Range("M62").Value = -2311
$ws.Range("H65").Value = 1981.9474
$ws.Range("I65").Value = 999
$ws.Range("K65").Value = 8991
$ws.Range("M65").Value = -5559
$ws.Range("H68").Value = 1170.2727
$ws.Range("I68").Value = 1134
$ws.Range("J68").Value = 1200.5
$ws.Range("K68").Value = 3402
$ws.Range("L68").Value = 3601.5
$ws.Range("M68").Value = -2591
$ws.Range("N68").Value = -5223.5
$ws.Range("H71").Value = 1170.2727
$ws.Range("I71").Value = 1134
$ws.Range("J71").Value = 1200.5
$ws.Range("K71").Value = 10206
$ws.Range("L71").Value = 10804.5
$ws.Range("M71").Value = -6150
$ws.Range("N71").Value = -18916.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("K43").Value = 1000
$ws.Range("M43").Value = -849
$ws.Range("H70").Value = 599
$ws.Range("I70").Value = 599
$ws.Range("K70").Value = 599
$ws.Range("M70").Value = -329
$ws.Range("H73").Value = 599
$ws.Range("I73").Value = 599
$ws.Range("K73").Value = 599
$ws.Range("M73").Value = 337
$ws.Range("H80").Value = 3199.4443
$ws.Range("I80").Value = 2359.2
$ws.Range("J80").Value = 4249.75
$ws.Range("K80").Value = 2359.2
$ws.Range("L80").Value = 4249.75
$ws.Range("M80").Value = -1361.2
$ws.Range("N80").Value = -6245.75
$ws.Range("H83").Value = 3199.4443
$ws.Range("I83").Value = 2359.2
$ws.Range("J83").Value = 4249.75
$ws.Range("K83").Value = 11796
$ws.Range("L83").Value = 21248.75
$ws.Range("M83").Value = -6804
$ws.Range("N83").Value = -31232.75
$ws.Range("H100").Value = 22000
$ws.Range("J100").Value = 22000
$ws.Range("L100").Value = 22000
$ws.Range("N100").Value = -24164

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 21748918
$ws.Range("I40").Value = 22736188
$ws.Range("K40").Value = 22736188
$ws.Range("M40").Value = -22736052
$ws.Range("H46").Value = 4818.476
$ws.Range("I46").Value = 1373.75
$ws.Range("J46").Value = 5629
$ws.Range("K46").Value = 1373.75
$ws.Range("L46").Value = 5629
$ws.Range("M46").Value = -1185.75
$ws.Range("N46").Value = -6005
$ws.Range("H50").Value = 1500000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("H122").Value = 31254154
$ws.Range("I122").Value = 31254154
$ws.Range("K122").Value = 93762462
$ws.Range("M122").Value = -93760012
$ws.Range("H136").Value = 3387.111
$ws.Range("I136").Value = 2393.476
$ws.Range("K136").Value = 7180.428
$ws.Range("M136").Value = -4630.428

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 128087.25
$ws.Range("I2").Value = 2339.8
$ws.Range("K2").Value = 2339.8
$ws.Range("M2").Value = -2227.8
$ws.Range("H5").Value = 6309.057
$ws.Range("I5").Value = 3181.5789
$ws.Range("K5").Value = 3181.5789
$ws.Range("M5").Value = -3069.5789
$ws.Range("H18").Value = 69695.5
$ws.Range("J18").Value = 69695.5
$ws.Range("L18").Value = 69695.5
$ws.Range("N18").Value = -70041.5
$ws.Range("H100").Value = 1875.5
$ws.Range("I100").Value = 1586.7222
$ws.Range("K100").Value = 3173.4444
$ws.Range("M100").Value = -2632.4444
$ws.Range("H107").Value = 2571.4583
$ws.Range("I107").Value = 2465.8
$ws.Range("J107").Value = 3099.75
$ws.Range("K107").Value = 7397.400000000001
$ws.Range("L107").Value = 9299.25
$ws.Range("M107").Value = -5477.400000000001
$ws.Range("N107").Value = -13139.25
$ws.Range("H132").Value = 12824632
$ws.Range("I132").Value = 17859916
$ws.Range("J132").Value = 7545.364
$ws.Range("K132").Value = 53579748
$ws.Range("L132").Value = 22636.092
$ws.Range("M132").Value = -53577218
$ws.Range("N132").Value = -27696.092

# --- Remove cell LTW!N50 (column dropped from output row) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N50").ClearContents()
